$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44189
$ws.Range("K2").Value = "Dina"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16562
$ws.Range("Q2").Value = "`$/caja 18 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 920
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 44176
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17400
$ws.Range("Q3").Value = "`$/caja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 967
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 44537
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 21500
$ws.Range("P4").Value = 21250
$ws.Range("Q4").Value = "`$/caja 15 kilos"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1417
$ws.Range("T4").Value = 15

$ws.Range("D5").Value = 44174
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 75
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9467
$ws.Range("Q5").Value = "`$/caja 10 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 947
$ws.Range("T5").Value = 10

$ws.Range("D7").Value = 44187
$ws.Range("K7").Value = "Dina"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 15455
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1030
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44165
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 16500
$ws.Range("Q8").Value = "`$/caja 15 kilos granel"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 1100
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 44544
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 16000
$ws.Range("O9").Value = 17000
$ws.Range("P9").Value = 16500
$ws.Range("Q9").Value = "`$/caja 15 kilos"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1100
$ws.Range("T9").Value = 15

$ws.Range("D10").Value = 44551
$ws.Range("K10").Value = "Castle Brite"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 15500
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15750
$ws.Range("Q10").Value = "`$/caja 15 kilos"
$ws.Range("R10").Value = "Región de O'Higgins"
$ws.Range("S10").Value = 1050
$ws.Range("T10").Value = 15

$ws.Range("D11").Value = 44181
$ws.Range("K11").Value = "Modesto"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("Q11").Value = "`$/caja 18 kilos"
$ws.Range("R11").Value = "Región de Coquimbo"
$ws.Range("S11").Value = 1139
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44552
$ws.Range("K12").Value = "Castle Brite"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 15500
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15750
$ws.Range("Q12").Value = "`$/caja 15 kilos"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 1050
$ws.Range("T12").Value = 15
